$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Datos actualizados" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 19 de Mayo de 2020 a las 09:35"

# Row 30 - Singapur
$ws.Range("B30").Value = 28794
$ws.Range("C30").Value = 451
$ws.Range("E30").Value = 18937

# Row 37 - Rumania
$ws.Range("E37").Value = 5982
$ws.Range("G37").Value = 4
$ws.Range("H37").Value = 1124

# Row 66 - Armenia
$ws.Range("B66").Value = 5041
$ws.Range("C66").Value = 218
$ws.Range("D66").Value = 2164
$ws.Range("E66").Value = 2813
$ws.Range("G66").Value = 3
$ws.Range("H66").Value = 64

# Row 105 - Letonia
$ws.Range("B105").Value = 1012
$ws.Range("C105").Value = 3
$ws.Range("D105").Value = 694
$ws.Range("E105").Value = 297
$ws.Range("G105").Value = 2
$ws.Range("H105").Value = 21
